$d = $word.ActiveDocument

# Change 1: text replacement in the first paragraph
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    "La carrera espacial entre Estados Unidos y Rusia, el juego consiste",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El juego se basa en la carrera espacial entre Estados Unidos y Rusia, el juego consiste",
    2)
